$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s_vals data (regen to filter save games)
$data = @{
    2 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    3 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 645.3272768299601, 648.4758912050064)
    4 = @(0.6753301551942219, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 1.372039145084537)
    5 = @(1.459612070389937, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 2.156321060280252)
    6 = @(0.0001488876196638067, 0.3127903958511391, 3.900430680208489, 8.660232485948974, 12.87360244962827)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
